# "Fixed typo in image 3.3.7-indices.png #24"
#
# The authoritative content fix in this commit is a swap of two text
# values inside a small 2-column lookup table: the cell that read
# "2+4" should read "1+5", and the cell that read "1+5" should read
# "2+4". Everything else in the diff (the cached datetimeFigureOut
# field text on the slide master/layouts, and the random p14:modId
# "last modified" nonces PowerPoint stamps onto every graphicFrame on
# save) is incidental save-time metadata, not an intentional edit, so
# we still touch the date placeholders (the object model lets us) but
# we don't try to forge the random modId values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Fix the swapped "2+4" / "1+5" table cells -----------------------
foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cell = $tbl.Cell($r, $c)
                $tr = $cell.Shape.TextFrame.TextRange
                if ($tr.Text -eq "2+4") {
                    $tr.Text = "1+5"
                } elseif ($tr.Text -eq "1+5") {
                    $tr.Text = "2+4"
                }
            }
        }
    }
}

# --- 2. Refresh the cached "today" date shown in the Date placeholders --
# (slide master + every custom layout use the same cached literal text)
$newDate = "11/27/24"

$sm = $p.SlideMaster
foreach ($shp in $sm.Shapes) {
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $sm.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    foreach ($shp in $cl.Shapes) {
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
